$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$colA = 'products__item\ in-stock\ products__item_3-in-row"]:nth-child(2) [type="button'
$colB = 'New!iPhone 15 Silicone Case'
$colC = 'Buy'
$colD = [char]0x5C + "31 52184-case-685"
$colE = [char]0x5C + "31 52185-case-689"
$colF = 'iPhone 15 Plus'
$colG = 'Cypress'

foreach ($r in 6,7) {
    $ws.Cells.Item($r, 1).Value = $colA
    $ws.Cells.Item($r, 2).Value = $colB
    $ws.Cells.Item($r, 3).Value = $colC
    $ws.Cells.Item($r, 4).Value = $colD
    $ws.Cells.Item($r, 5).Value = $colE
    $ws.Cells.Item($r, 6).Value = $colF
    $ws.Cells.Item($r, 7).Value = $colG
}
